$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add the new "Informe Final de SQA" document entry for row 25 (E205)
$ws.Range("B25").Value = "Informe Final de SQA"

# Update the active window's scroll position and selection to match the edit
$excel.ActiveWindow.ScrollRow = 16
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D25").Select()
